# Fix MIN/MAX ("/") placeholder values on the DEF sheet for several
# DEF_DEF rows: they should contain real numeric MIN/MAX bounds instead
# of the "/" string placeholder (shared string index 46).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DEF")

# row -> (MIN, MAX)
$updates = @{
    2  = @(0, 0)
    4  = @(0, 0)
    6  = @(0, 0)
    8  = @(0, 0)
    10 = @(0, 1000000)
    12 = @(1, 128)
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # column D = MIN
    $ws.Cells.Item($row, 5).Value = $vals[1]   # column E = MAX
}
